$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.684.65"
$ws.Range("E2").Value = "  -0.96%  "

$ws.Range("D3").Value = "2.538.07"
$ws.Range("E3").Value = "  -1.04%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "308.62"
$ws.Range("E5").Value = "  -2.02%  "

$ws.Range("D6").Value = "100.83"
$ws.Range("E6").Value = "  +4.34%  "

$ws.Range("D7").Value = "0.571"
$ws.Range("E7").Value = "  -1.18%  "

$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("E9").Value = "  -1.83%  "

$ws.Range("D10").Value = "36.01"
$ws.Range("E10").Value = "  +1.67%  "

$ws.Range("E11").Value = "  -1.02%  "

$ws.Range("D12").Value = "7.37"
$ws.Range("E12").Value = "  -0.95%  "

$ws.Range("E13").Value = "  +0.11%  "

$ws.Range("D14").Value = "2.932.05"
$ws.Range("E14").Value = "  -0.90%  "

$ws.Range("D15").Value = "15.92"
$ws.Range("E15").Value = "  +5.79%  "

$ws.Range("D16").Value = "2.491.41"
$ws.Range("E16").Value = "  -5.28%  "

$ws.Range("D17").Value = "0.815"
$ws.Range("E17").Value = "  -3.15%  "

$ws.Range("D18").Value = "42.688.16"
$ws.Range("E18").Value = "  -1.00%  "

$ws.Range("D19").Value = "6.78"
$ws.Range("E19").Value = "  -0.61%  "

$ws.Range("D20").Value = "0.0₃0954"
$ws.Range("E20").Value = "  -0.68%  "

$ws.Range("D21").Value = "12.23"
$ws.Range("E21").Value = "  -2.61%  "

$ws.Range("D22").Value = "69.52"
$ws.Range("E22").Value = "  +0.40%  "

$ws.Range("D23").Value = "244.07"
$ws.Range("E23").Value = "  -3.63%  "

$ws.Range("D24").Value = "2.89"
$ws.Range("E24").Value = "  -2.13%  "

$ws.Range("D25").Value = "2.05"
$ws.Range("E25").Value = "  -1.25%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("D27").Value = "26.09"
$ws.Range("E27").Value = "  -2.73%  "

$ws.Range("E28").Value = "  -3.80%  "

$ws.Range("D29").Value = "39.32"
$ws.Range("E29").Value = "  -1.63%  "

$ws.Range("D30").Value = "10.16"
$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("D31").Value = "5.80"
$ws.Range("E31").Value = "  -0.47%  "

$ws.Range("D32").Value = "156.49"
$ws.Range("E32").Value = "  +0.97%  "

$ws.Range("D33").Value = "2.74"
$ws.Range("E33").Value = "  +11.99%  "

$ws.Range("D34").Value = "0.0794"
$ws.Range("E34").Value = "  -1.59%  "

$ws.Range("E35").Value = "  -3.09%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "2.03"
$ws.Range("E36").Value = "  -4.21%  "

$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").Value = "18.35"
$ws.Range("E37").Value = "  -3.06%  "

$ws.Range("D38").Value = "3.17"
$ws.Range("E38").Value = "  -6.50%  "

$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("E40").Value = "  +0.80%  "

$ws.Range("D41").Value = "4.33"
$ws.Range("E41").Value = "  +8.65%  "

$ws.Range("D42").Value = "21.92"
$ws.Range("E42").Value = "  -2.71%  "

$ws.Range("E43").Value = "  +0.09%  "

$ws.Range("E44").Value = "  +1.87%  "

$ws.Range("D45").Value = "0.0299"
$ws.Range("E45").Value = "  -1.76%  "

$ws.Range("D46").Value = "1.973.62"
$ws.Range("E46").Value = "  -1.64%  "

$ws.Range("D47").Value = "8.87"
$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("D48").Value = "81.17"
$ws.Range("E48").Value = "  -1.79%  "

$ws.Range("D49").Value = "0.193"
$ws.Range("E49").Value = "  -0.33%  "

$ws.Range("E50").Value = "  +10.90%  "

$ws.Range("D51").Value = "2.727.08"
$ws.Range("E51").Value = "  -2.97%  "

